$d = $word.ActiveDocument

# 1. Update the two "date de visite" table cells: 14/06/2018 -> 15/06/2018
#    (wdReplaceAll replaces every occurrence in the story, i.e. both cells)
$d.Content.Find.Execute("14/06/2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "15/06/2018", 2)

# 2. Update the "Montant total" amount in the recap table: -8 -> 10405,00
$d.Content.Find.Execute("-8", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "10405,00", 2)

# 3. Update the final "Montant total" figure at the end of the document:
#    51007,00 -> 61420,00
$d.Content.Find.Execute("51007,00", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "61420,00", 2)
